$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 20002196
$ws.Range("I86").Value = 38463616
$ws.Range("J86").Value = 2324.0833
$ws.Range("K86").Value = 38463616
$ws.Range("L86").Value = 2324.0833
$ws.Range("M86").Value = -38462493
$ws.Range("N86").Value = -4570.0833
$ws.Range("H89").Value = 20002196
$ws.Range("I89").Value = 38463616
$ws.Range("J89").Value = 2324.0833
$ws.Range("K89").Value = 192318080
$ws.Range("L89").Value = 11620.4165
$ws.Range("M89").Value = -192312464
$ws.Range("N89").Value = -22852.4165
$ws.Range("H132").Value = 3953.4182
$ws.Range("I132").Value = 4300.023
$ws.Range("J132").Value = 2567
$ws.Range("K132").Value = 12900.069
$ws.Range("L132").Value = 7701
$ws.Range("M132").Value = -10370.069
$ws.Range("N132").Value = -12761
$ws.Range("H137").Value = 10526934
$ws.Range("I137").Value = 603.6923
$ws.Range("J137").Value = 33333984
$ws.Range("K137").Value = 1811.0769
$ws.Range("L137").Value = 100001952
$ws.Range("M137").Value = 738.9231
$ws.Range("N137").Value = -100007052

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").ClearContents()
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = 0
$ws.Range("H36").Value = 3545.2
$ws.Range("I36").Value = 3545.2
$ws.Range("K36").Value = 3545.2
$ws.Range("M36").Value = -3199.2
$ws.Range("H45").Value = 1526.6552
$ws.Range("I45").Value = 1529.9412
$ws.Range("J45").Value = 1522
$ws.Range("K45").Value = 1529.9412
$ws.Range("L45").Value = 1522
$ws.Range("M45").Value = -1152.9412
$ws.Range("N45").Value = -2276
$ws.Range("H61").Value = 6173886.5
$ws.Range("I61").Value = 7576662.5
$ws.Range("J61").Value = 1673.5333
$ws.Range("K61").Value = 7576662.5
$ws.Range("L61").Value = 1673.5333
$ws.Range("M61").Value = -7576450.5
$ws.Range("N61").Value = -2097.5333
$ws.Range("H132").Value = 5436872.5
$ws.Range("I132").Value = 6758520.5
$ws.Range("J132").Value = 3429.7778
$ws.Range("K132").Value = 20275561.5
$ws.Range("L132").Value = 10289.3334
$ws.Range("M132").Value = -20273031.5
$ws.Range("N132").Value = -15349.3334
$ws.Range("H136").Value = 6173886.5
$ws.Range("I136").Value = 7576662.5
$ws.Range("J136").Value = 1673.5333
$ws.Range("K136").Value = 22729987.5
$ws.Range("L136").Value = 5020.5999
$ws.Range("M136").Value = -22727437.5
$ws.Range("N136").Value = -10120.5999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H118").Value = 7906.923
$ws.Range("J118").Value = 7906.923
$ws.Range("L118").Value = 7906.923
$ws.Range("N118").Value = -11220.923
$ws.Range("H119").Value = 35341.645
$ws.Range("J119").Value = 35341.645
$ws.Range("L119").Value = 35341.645
$ws.Range("N119").Value = -45017.645
$ws.Range("H134").Value = 4578.794
$ws.Range("I134").Value = 3394.04
$ws.Range("K134").Value = 10182.12
$ws.Range("M134").Value = -7647.119999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1052012
$ws.Range("I6").Value = 1052012
$ws.Range("K6").Value = 1052012
$ws.Range("M6").Value = -1051899
$ws.Range("H31").Value = 6539992
$ws.Range("I31").Value = 4029.976
$ws.Range("J31").Value = 37041148
$ws.Range("K31").Value = 4029.976
$ws.Range("L31").Value = 37041148
$ws.Range("M31").Value = -3734.976
$ws.Range("N31").Value = -37041738
$ws.Range("H34").Value = 6539992
$ws.Range("I34").Value = 4029.976
$ws.Range("J34").Value = 37041148
$ws.Range("K34").Value = 4029.976
$ws.Range("L34").Value = 37041148
$ws.Range("M34").Value = -3827.976
$ws.Range("N34").Value = -37041552
$ws.Range("H99").Value = 1106.2
$ws.Range("I99").Value = 1318.6666
$ws.Range("J99").Value = 787.5
$ws.Range("K99").Value = 1318.6666
$ws.Range("L99").Value = 787.5
$ws.Range("M99").Value = 179.3334
$ws.Range("N99").Value = -3783.5
$ws.Range("H126").Value = 1106.2
$ws.Range("I126").Value = 1318.6666
$ws.Range("J126").Value = 787.5
$ws.Range("K126").Value = 3955.9998
$ws.Range("L126").Value = 2362.5
$ws.Range("M126").Value = -1485.9998
$ws.Range("N126").Value = -7302.5
$ws.Range("H132").Value = 8930309
$ws.Range("I132").Value = 11906179
$ws.Range("J132").Value = 2697.2856
$ws.Range("K132").Value = 35718537
$ws.Range("L132").Value = 8091.8568
$ws.Range("M132").Value = -35716007
$ws.Range("N132").Value = -13151.8568
$ws.Range("H134").Value = 1661.1
$ws.Range("I134").Value = 1588
$ws.Range("J134").Value = 1953.5
$ws.Range("K134").Value = 4764
$ws.Range("L134").Value = 5860.5
$ws.Range("M134").Value = -2229
$ws.Range("N134").Value = -10930.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3019.9
$ws.Range("I5").Value = 298.5
$ws.Range("K5").Value = 895.5
$ws.Range("M5").Value = -783.5
$ws.Range("H113").Value = 1385.091
$ws.Range("I113").Value = 463.1111
$ws.Range("J113").Value = 2023.3846
$ws.Range("K113").Value = 1389.3333
$ws.Range("L113").Value = 6070.1538
$ws.Range("M113").Value = 780.6667
$ws.Range("N113").Value = -10410.1538
$ws.Range("H135").Value = 3019.9
$ws.Range("I135").Value = 298.5
$ws.Range("K135").Value = 2686.5
$ws.Range("M135").Value = -151.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 23690
$ws.Range("I70").Value = 43360
$ws.Range("J70").Value = 4020
$ws.Range("K70").Value = 43360
$ws.Range("L70").Value = 4020
$ws.Range("M70").Value = -43090
$ws.Range("N70").Value = -4560
$ws.Range("H73").Value = 23690
$ws.Range("I73").Value = 43360
$ws.Range("J73").Value = 4020
$ws.Range("K73").Value = 43360
$ws.Range("L73").Value = 4020
$ws.Range("M73").Value = -42424
$ws.Range("N73").Value = -5892
$ws.Range("H126").Value = 5745.4546
$ws.Range("I126").Value = 4333.3335
$ws.Range("J126").Value = 6275
$ws.Range("K126").Value = 13000.0005
$ws.Range("L126").Value = 18825
$ws.Range("M126").Value = -10530.0005
$ws.Range("N126").Value = -23765
$ws.Range("H132").Value = 4208.732
$ws.Range("I132").Value = 3198
$ws.Range("J132").Value = 6158
$ws.Range("K132").Value = 9594
$ws.Range("L132").Value = 18474
$ws.Range("M132").Value = -7064
$ws.Range("N132").Value = -23534

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4433.846
$ws.Range("I122").Value = 4123.2256
$ws.Range("J122").Value = 5637.5
$ws.Range("K122").Value = 12369.6768
$ws.Range("L122").Value = 16912.5
$ws.Range("M122").Value = -9919.676799999999
$ws.Range("N122").Value = -21812.5
$ws.Range("H136").Value = 27785108
$ws.Range("I136").Value = 35716070
$ws.Range("K136").Value = 107148210
$ws.Range("M136").Value = -107145660

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 5525.5
$ws.Range("I61").Value = 1051
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 1051
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -759
$ws.Range("N61").Value = -10584
